# USDM CT deliverable sprint rollover: rename the "Sp12" terminology-change
# tabs to "Sp13" (add/remove/change tracking sheets are rolled forward each
# sprint; the newest pair keeps the generic "new"/"chg" sheetId/rId slots and
# simply gets relabeled for the current sprint number).

$wb = $excel.ActiveWorkbook

$newSheet = $wb.Worksheets.Item("Terminology Changes Sp12 - new")
$newSheet.Name = "Terminology Changes Sp13 - new"

$chgSheet = $wb.Worksheets.Item("Terminology Changes Sp12 - chg")
$chgSheet.Name = "Terminology Changes Sp13 - chg"
